# Refresh crypto price/volume data in the active worksheet.
# Column layout: A=index (unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Coin, Link, Price, Volume(1h)
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '23.017.17', '  -3.77%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.601.33', '  -2.72%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  -0.01%  ')
    ,@(5, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.06%  ')
    ,@(6, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '301.16', '  -2.95%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.3775', '  -2.81%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3632', '  -5.13%  ')
    ,@(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '49.84', '  -1.42%  ')
    ,@(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.258', '  -5.55%  ')
    ,@(11, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.01%  ')
    ,@(12, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08114', '  -3.50%  ')
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '22.59', '  -5.10%  ')
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.590', '  -5.76%  ')
    ,@(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.328', '  -6.66%  ')
    ,@(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001242', '  -5.41%  ')
    ,@(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.604.39', '  -2.41%  ')
    ,@(18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '91.95', '  -2.05%  ')
    ,@(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06816', '  -1.99%  ')
    ,@(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '18.22', '  -6.70%  ')
    ,@(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.554', '  -4.89%  ')
    ,@(22, 'BitDAO', 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit', '0.5577', '  -5.44%  ')
    ,@(23, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.001', '  +0.04%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '13.14', '  -3.61%  ')
    ,@(25, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '23.014.12', '  -3.81%  ')
    ,@(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.364', '  -2.82%  ')
    ,@(27, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.829', '  -2.60%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '21.06', '  -3.72%  ')
    ,@(29, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '150.92', '  -1.60%  ')
    ,@(30, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '5.247', '  -4.84%  ')
    ,@(31, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '133.85', '  -2.04%  ')
    ,@(32, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.236', '  -10.52%  ')
    ,@(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.811', '  -10.67%  ')
    ,@(34, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.777.07', '  -2.67%  ')
    ,@(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.9670', '  -1.53%  ')
    ,@(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.07577', '  -5.70%  ')
    ,@(37, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '10.41', '  -0.04%  ')
    ,@(38, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '6.267', '  -4.74%  ')
    ,@(39, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2532', '  -5.18%  ')
    ,@(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02698', '  -7.39%  ')
    ,@(41, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08878', '  -2.32%  ')
    ,@(42, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.367', '  -3.31%  ')
    ,@(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.7006', '  -6.97%  ')
    ,@(44, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '12.40', '  -6.99%  ')
    ,@(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '15.28', '  -7.11%  ')
    ,@(46, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6608', '  -4.45%  ')
    ,@(47, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9991', '  -0.15%  ')
    ,@(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.302', '  -4.98%  ')
    ,@(49, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.987', '  -2.54%  ')
    ,@(50, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '132.14', '  -1.83%  ')
    ,@(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07905', '  -4.28%  ')
)

foreach ($item in $data) {
    $row = $item[0]
    $coin = $item[1]
    $link = $item[2]
    $price = $item[3]
    $volume = $item[4]

    $ws.Cells.Item($row, 2).Value = $coin
    $ws.Cells.Item($row, 3).Value = $link

    # Force text format on the Price cell so numeric-looking strings
    # (e.g. '301.16', '1.001', '0.9670') are preserved verbatim instead
    # of being re-interpreted/rounded by Excel as a number. Values such as
    # '23.017.17' contain two dots and Excel already keeps those as text.
    if ($price -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Cells.Item($row, 4).NumberFormat = '@'
    }
    $ws.Cells.Item($row, 4).Value = $price

    $ws.Cells.Item($row, 5).Value = $volume
}
